# Auto-generated Excel COM-interop script to update the cryptos list
# per the commit "Updated cryptos list on Mon Apr 17 19:23:52 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.703.32'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '2.097.56'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '344.48'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = '0.5137'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").Value = '0.4408'
$ws.Range("E8").Value = '  -3.59%  '
$ws.Range("D9").Value = '52.68'
$ws.Range("E9").Value = '  -2.97%  '
$ws.Range("D10").Value = '0.09182'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = '1.172'
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").Value = '24.97'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '2.098.99'
$ws.Range("E13").Value = '  -2.11%  '
$ws.Range("D14").Value = '8.287'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").Value = '6.757'
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("D16").Value = '99.66'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = '1.008'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '20.85'
$ws.Range("E19").Value = '  +6.22%  '
$ws.Range("D20").Value = '0.06624'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '1.006'
$ws.Range("D22").Value = '6.197'
$ws.Range("E22").Value = '  -3.21%  '
$ws.Range("D23").Value = '29.750.24'
$ws.Range("E23").Value = '  -3.49%  '
$ws.Range("D24").Value = '12.62'
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").Value = '2.318'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").Value = '2.347.89'
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").Value = '21.90'
$ws.Range("E27").Value = '  -3.55%  '
$ws.Range("D28").Value = '2.531'
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("D29").Value = '161.96'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '132.99'
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").Value = '1.129'
$ws.Range("E31").Value = '  -7.77%  '
$ws.Range("E32").Value = '  -3.31%  '
$ws.Range("D33").Value = '1.653'
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("D34").Value = '6.172'
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("E35").Value = '  -2.04%  '
$ws.Range("D36").Value = '10.45'
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = '6.020'
$ws.Range("E37").Value = '  -2.40%  '
$ws.Range("D38").Value = '0.02569'
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("D39").Value = '0.06728'
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").Value = '12.45'
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").Value = '0.6880'
$ws.Range("D42").Value = '0.2237'
$ws.Range("E42").Value = '  -4.50%  '
$ws.Range("D43").Value = '1.297'
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").Value = '0.6657'
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("D45").Value = '14.30'
$ws.Range("E45").Value = '  -3.65%  '
$ws.Range("D46").Value = '2.306'
$ws.Range("E46").Value = '  -2.56%  '
$ws.Range("D47").Value = '3.610'
$ws.Range("E47").Value = '  -4.00%  '
$ws.Range("D48").Value = '0.00000000349'
$ws.Range("E48").Value = '  -5.89%  '
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").Value = '0.3376'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '82.21'
$ws.Range("E51").Value = '  -1.54%  '
